# "Modification in login page"
#
# The "login" sheet contained a duplicate row (A2/B2 = 9876543211 / admin).
# That row is removed, so every row below it shifts up by one; the row that
# used to be the last one (16) becomes blank, and the hyperlink that used to
# sit on the "987654321@" row (old row 8) now lives on row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# Drop the hyperlink before we move things around so it doesn't keep
# pointing at a now-shifted cell.
$ws.Range("A8").Hyperlinks.Delete()

# Remove row 2 entirely; rows 3-16 shift up to become rows 2-15.
$ws.Rows.Item(2).Delete()

# Former row 16 is now empty, but it keeps column A's left-aligned format,
# same as the other data rows.
$ws.Range("A16").HorizontalAlignment = -4131

# Re-create the hyperlink on its new location (former row 8 -> row 7) and
# restore the built-in "Hyperlink" look (blue, underlined).
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:987654321@")
$ws.Range("A7").Style = "Hyperlink"

# Match the captured selection state of the sheet after the edit.
$ws.Range("A16").Select()
